# Add a "payment type" row to the customer card sheet.
# The sheet has a repeating "label (merged A:B) / value (merged C:F)" row
# pattern (rows 4, 6, 8, 11, 13, 15 ...). Row 14 is currently an empty
# spacer row between "Плательщик" (row 13) and "Дата заказа" (row 15);
# we turn it into a new "Тип оплаты" / "@paymentType" row, matching the
# look of the other label/value rows exactly by cloning row 8's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (font, fill, borders, alignment, number format,
# row height) of an existing label/value row onto row 14.
$ws.Range("A8:F8").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the merges for the label cell and the value cell, same as
# the other rows of this pattern.
$ws.Range("A14:B14").Merge()
$ws.Range("C14:F14").Merge()

# Fill in the new labels.
$ws.Range("A14").Value = "Тип оплаты"
$ws.Range("C14").Value = "@paymentType"
